$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Slovakia" sheet by copying "Portugal" (it already
#    has the narrower column widths that Slovakia needs) and placing
#    the copy after the last sheet.
# ------------------------------------------------------------------
$portugal = $wb.Worksheets.Item("Portugal")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Copy($null, $lastSheet)

$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# ------------------------------------------------------------------
# 2. Slovakia needs an extra row (like Germany/Czech: FBI800, PCH800,
#    Wg, Slot Cards) instead of Portugal's 3-row tail (PCH800, Wg,
#    Slot Cards). Insert a row above the old "PCH800" row so row 8
#    becomes free, and copy that row's formatting onto it.
# ------------------------------------------------------------------
$slovakia.Rows.Item(8).Insert()
[void]$slovakia.Range("A9").Copy()
[void]$slovakia.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

# Rows 3-5 inherited Portugal's custom (wrapped-text) row height; with
# Slovakia's own content the text no longer needs the extra height, so
# restore the default row height by auto-fitting.
[void]$slovakia.Rows.Item(3).AutoFit()
[void]$slovakia.Rows.Item(4).AutoFit()
[void]$slovakia.Rows.Item(5).AutoFit()

# ------------------------------------------------------------------
# 3. Fill in the Slovakia-specific content.
# ------------------------------------------------------------------
$slovakia.Range("A8").Value = "FBI800"
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3218"

# ------------------------------------------------------------------
# 4. Page setup (portrait) like the other market sheets already have.
# ------------------------------------------------------------------
$portugal.PageSetup.Orientation = 1
$slovakia.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 5. View state.
#    - Germany's selection moved from A10 to A8:A11.
#    - Portugal is no longer the active tab; its selection moves to B2.
#    - Slovakia becomes the active tab, selection A8:A11 (active cell A8).
# ------------------------------------------------------------------
$germany = $wb.Worksheets.Item("Germany")
[void]$germany.Select()
[void]$germany.Range("A8:A11").Select()

[void]$portugal.Select()
[void]$portugal.Range("B2").Select()

[void]$slovakia.Select()
[void]$slovakia.Range("A8:A11").Select()
